$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 140 - this shifts the existing rows 140:169 down to 141:170
# (Excel also grows the used range / dimension automatically).
$ws.Rows.Item(140).Insert()

# Populate the new row 140 with the latest weekly price record.
$ws.Range("A140").Value = 7
$ws.Range("B140").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C140").Value = 'Ñuble'
$ws.Range("D140").Value = 45209
$ws.Range("E140").Value = 16
$ws.Range("F140").Value = 100112031
$ws.Range("G140").Value = 'Poroto verde'
$ws.Range("H140").Value = 'Sin especificar'
$ws.Range("I140").Value = 'Primera'
$ws.Range("J140").Value = 30
$ws.Range("K140").Value = 28000
$ws.Range("L140").Value = 28000
$ws.Range("M140").Value = 28000
$ws.Range("N140").Value = '$/malla 25 kilos'
$ws.Range("O140").Value = 'Perú'
$ws.Range("P140").Value = 1120
$ws.Range("Q140").Value = 25
$ws.Range("R140").Value = 'Hortaliza'
